# Scheduled-runner style refresh of the per-sheet Leve profit data
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns, H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
#
# Columns: H=8 I=9 J=10 K=11 L=12 M=13 N=14
# Setting a cell's Value to "" clears it so the saved XML drops the <c>
# element entirely (matches rows where the profit column has no cached
# value for a given item).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(62, 8).Value = 15364.771   # H62: 8546.052 -> 15364.771
$ws.Cells.Item(62, 9).Value = 16718.291   # I62: 9163.82 -> 16718.291
$ws.Cells.Item(62, 10).Value = 4875   # J62: 4685 -> 4875
$ws.Cells.Item(62, 11).Value = 16718.291   # K62: 9163.82 -> 16718.291
$ws.Cells.Item(62, 12).Value = 4875   # L62: 4685 -> 4875
$ws.Cells.Item(62, 13).Value = -16094.291   # M62: -8539.82 -> -16094.291
$ws.Cells.Item(62, 14).Value = -6123   # N62: -5933 -> -6123

$ws.Cells.Item(65, 8).Value = 15364.771   # H65: 8546.052 -> 15364.771
$ws.Cells.Item(65, 9).Value = 16718.291   # I65: 9163.82 -> 16718.291
$ws.Cells.Item(65, 10).Value = 4875   # J65: 4685 -> 4875
$ws.Cells.Item(65, 11).Value = 83591.455   # K65: 45819.1 -> 83591.455
$ws.Cells.Item(65, 12).Value = 24375   # L65: 23425 -> 24375
$ws.Cells.Item(65, 13).Value = -80471.455   # M65: -42699.1 -> -80471.455
$ws.Cells.Item(65, 14).Value = -30615   # N65: -29665 -> -30615

$ws.Cells.Item(136, 8).Value = 60310   # H136: 60486.668 -> 60310
$ws.Cells.Item(136, 10).Value = 60310   # J136: 60486.668 -> 60310
$ws.Cells.Item(136, 12).Value = 60310   # L136: 60486.668 -> 60310
$ws.Cells.Item(136, 14).Value = -70510   # N136: -70686.66800000001 -> -70510

$ws.Cells.Item(138, 8).Value = 1992.21   # H138: 1973.7677 -> 1992.21
$ws.Cells.Item(138, 9).Value = 1148.3334   # I138: 1108.3784 -> 1148.3334
$ws.Cells.Item(138, 10).Value = 2466.8906   # J138: 2490.2097 -> 2466.8906
$ws.Cells.Item(138, 11).Value = 3445.0002   # K138: 3325.1352 -> 3445.0002
$ws.Cells.Item(138, 12).Value = 7400.6718   # L138: 7470.6291 -> 7400.6718
$ws.Cells.Item(138, 13).Value = 1694.9998   # M138: 1814.8648 -> 1694.9998
$ws.Cells.Item(138, 14).Value = -17680.6718   # N138: -17750.6291 -> -17680.6718

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(122, 8).Value = 1041.619   # H122: 1162.4286 -> 1041.619
$ws.Cells.Item(122, 9).Value = 994.5714   # I122: 1189.1428 -> 994.5714
$ws.Cells.Item(122, 11).Value = 2983.7142   # K122: 3567.4284 -> 2983.7142
$ws.Cells.Item(122, 13).Value = -533.7142000000003   # M122: -1117.4284 -> -533.7142000000003

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 2950.5625   # H20: 3060.6 -> 2950.5625
$ws.Cells.Item(20, 9).Value = 3128.5715   # I20: 3142.8572 -> 3128.5715
$ws.Cells.Item(20, 10).Value = 1704.5   # J20: 1909 -> 1704.5
$ws.Cells.Item(20, 11).Value = 3128.5715   # K20: 3142.8572 -> 3128.5715
$ws.Cells.Item(20, 12).Value = 1704.5   # L20: 1909 -> 1704.5
$ws.Cells.Item(20, 13).Value = -2881.5715   # M20: -2895.8572 -> -2881.5715
$ws.Cells.Item(20, 14).Value = -2198.5   # N20: -2403 -> -2198.5

$ws.Cells.Item(68, 8).Value = 19950   # H68: 0 -> 19950
$ws.Cells.Item(68, 10).Value = 19950   # J68: 0 -> 19950
$ws.Cells.Item(68, 12).Value = 19950   # L68: 0 -> 19950
$ws.Cells.Item(68, 14).Value = -21572   # N68: None -> -21572

$ws.Cells.Item(69, 8).Value = 0   # H69: 20000 -> 0
$ws.Cells.Item(69, 10).Value = 0   # J69: 20000 -> 0
$ws.Cells.Item(69, 12).Value = 0   # L69: 20000 -> 0
$ws.Cells.Item(69, 14).Value = ""   # N69: remove (was -21622)

$ws.Cells.Item(71, 8).Value = 19950   # H71: 0 -> 19950
$ws.Cells.Item(71, 10).Value = 19950   # J71: 0 -> 19950
$ws.Cells.Item(71, 12).Value = 59850   # L71: 0 -> 59850
$ws.Cells.Item(71, 14).Value = -67962   # N71: None -> -67962

$ws.Cells.Item(72, 8).Value = 0   # H72: 20000 -> 0
$ws.Cells.Item(72, 10).Value = 0   # J72: 20000 -> 0
$ws.Cells.Item(72, 12).Value = 0   # L72: 60000 -> 0
$ws.Cells.Item(72, 14).Value = ""   # N72: remove (was -68112)

$ws.Cells.Item(82, 8).Value = 12448.917   # H82: 10045.467 -> 12448.917
$ws.Cells.Item(82, 9).Value = 5430.1113   # I82: 3486.0833 -> 5430.1113
$ws.Cells.Item(82, 10).Value = 33505.332   # J82: 36283 -> 33505.332
$ws.Cells.Item(82, 11).Value = 5430.1113   # K82: 3486.0833 -> 5430.1113
$ws.Cells.Item(82, 12).Value = 33505.332   # L82: 36283 -> 33505.332
$ws.Cells.Item(82, 13).Value = -5047.1113   # M82: -3103.0833 -> -5047.1113
$ws.Cells.Item(82, 14).Value = -34271.332   # N82: -37049 -> -34271.332

$ws.Cells.Item(85, 8).Value = 12448.917   # H85: 10045.467 -> 12448.917
$ws.Cells.Item(85, 9).Value = 5430.1113   # I85: 3486.0833 -> 5430.1113
$ws.Cells.Item(85, 10).Value = 33505.332   # J85: 36283 -> 33505.332
$ws.Cells.Item(85, 11).Value = 5430.1113   # K85: 3486.0833 -> 5430.1113
$ws.Cells.Item(85, 12).Value = 33505.332   # L85: 36283 -> 33505.332
$ws.Cells.Item(85, 13).Value = -4104.1113   # M85: -2160.0833 -> -4104.1113
$ws.Cells.Item(85, 14).Value = -36157.332   # N85: -38935 -> -36157.332

$ws.Cells.Item(88, 8).Value = 25000   # H88: 17655.5 -> 25000
$ws.Cells.Item(88, 9).Value = 0   # I88: 5311 -> 0
$ws.Cells.Item(88, 10).Value = 25000   # J88: 30000 -> 25000
$ws.Cells.Item(88, 11).Value = 0   # K88: 5311 -> 0
$ws.Cells.Item(88, 12).Value = 25000   # L88: 30000 -> 25000
$ws.Cells.Item(88, 13).Value = ""   # M88: remove (was -4905)
$ws.Cells.Item(88, 14).Value = -25812   # N88: -30812 -> -25812

$ws.Cells.Item(91, 8).Value = 25000   # H91: 17655.5 -> 25000
$ws.Cells.Item(91, 9).Value = 0   # I91: 5311 -> 0
$ws.Cells.Item(91, 10).Value = 25000   # J91: 30000 -> 25000
$ws.Cells.Item(91, 11).Value = 0   # K91: 5311 -> 0
$ws.Cells.Item(91, 12).Value = 25000   # L91: 30000 -> 25000
$ws.Cells.Item(91, 13).Value = ""   # M91: remove (was -3907)
$ws.Cells.Item(91, 14).Value = -27808   # N91: -32808 -> -27808

$ws.Cells.Item(132, 8).Value = 500020000   # H132: 250036750 -> 500020000
$ws.Cells.Item(132, 10).Value = 500020000   # J132: 250036750 -> 500020000
$ws.Cells.Item(132, 12).Value = 500020000   # L132: 250036750 -> 500020000
$ws.Cells.Item(132, 14).Value = -500030120   # N132: -250046870 -> -500030120

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 6044.343   # H31: 5250.7607 -> 6044.343
$ws.Cells.Item(31, 9).Value = 1598.3636   # I31: 1374.1515 -> 1598.3636
$ws.Cells.Item(31, 10).Value = 13568.308   # J31: 15091.385 -> 13568.308
$ws.Cells.Item(31, 11).Value = 1598.3636   # K31: 1374.1515 -> 1598.3636
$ws.Cells.Item(31, 12).Value = 13568.308   # L31: 15091.385 -> 13568.308
$ws.Cells.Item(31, 13).Value = -1303.3636   # M31: -1079.1515 -> -1303.3636
$ws.Cells.Item(31, 14).Value = -14158.308   # N31: -15681.385 -> -14158.308

$ws.Cells.Item(34, 8).Value = 6044.343   # H34: 5250.7607 -> 6044.343
$ws.Cells.Item(34, 9).Value = 1598.3636   # I34: 1374.1515 -> 1598.3636
$ws.Cells.Item(34, 10).Value = 13568.308   # J34: 15091.385 -> 13568.308
$ws.Cells.Item(34, 11).Value = 1598.3636   # K34: 1374.1515 -> 1598.3636
$ws.Cells.Item(34, 12).Value = 13568.308   # L34: 15091.385 -> 13568.308
$ws.Cells.Item(34, 13).Value = -1396.3636   # M34: -1172.1515 -> -1396.3636
$ws.Cells.Item(34, 14).Value = -13972.308   # N34: -15495.385 -> -13972.308

$ws.Cells.Item(86, 8).Value = 31258510   # H86: 111128664 -> 31258510
$ws.Cells.Item(86, 9).Value = 50010360   # I86: 111128664 -> 50010360
$ws.Cells.Item(86, 10).Value = 5425   # J86: 0 -> 5425
$ws.Cells.Item(86, 11).Value = 50010360   # K86: 111128664 -> 50010360
$ws.Cells.Item(86, 12).Value = 5425   # L86: 0 -> 5425
$ws.Cells.Item(86, 13).Value = -50009237   # M86: -111127541 -> -50009237
$ws.Cells.Item(86, 14).Value = -7671   # N86: None -> -7671

$ws.Cells.Item(89, 8).Value = 31258510   # H89: 111128664 -> 31258510
$ws.Cells.Item(89, 9).Value = 50010360   # I89: 111128664 -> 50010360
$ws.Cells.Item(89, 10).Value = 5425   # J89: 0 -> 5425
$ws.Cells.Item(89, 11).Value = 250051800   # K89: 555643320 -> 250051800
$ws.Cells.Item(89, 12).Value = 27125   # L89: 0 -> 27125
$ws.Cells.Item(89, 13).Value = -250046184   # M89: -555637704 -> -250046184
$ws.Cells.Item(89, 14).Value = -38357   # N89: None -> -38357

$ws.Cells.Item(122, 8).Value = 1333.3334   # H122: 1014.2857 -> 1333.3334
$ws.Cells.Item(122, 9).Value = 1000   # I122: 1016.6667 -> 1000
$ws.Cells.Item(122, 10).Value = 2000   # J122: 1000 -> 2000
$ws.Cells.Item(122, 11).Value = 3000   # K122: 3050.0001 -> 3000
$ws.Cells.Item(122, 12).Value = 6000   # L122: 3000 -> 6000
$ws.Cells.Item(122, 13).Value = -550   # M122: -600.0001000000002 -> -550
$ws.Cells.Item(122, 14).Value = -10900   # N122: -7900 -> -10900

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(64, 8).Value = 1897422.1   # H64: 1597639.2 -> 1897422.1
$ws.Cells.Item(64, 9).Value = 1804.8   # I64: 988.875 -> 1804.8
$ws.Cells.Item(64, 10).Value = 2759066.2   # J64: 2758839.5 -> 2759066.2
$ws.Cells.Item(64, 11).Value = 5414.4   # K64: 2966.625 -> 5414.4
$ws.Cells.Item(64, 12).Value = 8277198.600000001   # L64: 8276518.5 -> 8277198.600000001
$ws.Cells.Item(64, 13).Value = -5144.4   # M64: -2696.625 -> -5144.4
$ws.Cells.Item(64, 14).Value = -8277738.600000001   # N64: -8277058.5 -> -8277738.600000001

$ws.Cells.Item(67, 8).Value = 1897422.1   # H67: 1597639.2 -> 1897422.1
$ws.Cells.Item(67, 9).Value = 1804.8   # I67: 988.875 -> 1804.8
$ws.Cells.Item(67, 10).Value = 2759066.2   # J67: 2758839.5 -> 2759066.2
$ws.Cells.Item(67, 11).Value = 5414.4   # K67: 2966.625 -> 5414.4
$ws.Cells.Item(67, 12).Value = 8277198.600000001   # L67: 8276518.5 -> 8277198.600000001
$ws.Cells.Item(67, 13).Value = -4478.4   # M67: -2030.625 -> -4478.4
$ws.Cells.Item(67, 14).Value = -8279070.600000001   # N67: -8278390.5 -> -8279070.600000001

$ws.Cells.Item(117, 8).Value = 1377   # H117: 1477.5385 -> 1377
$ws.Cells.Item(117, 9).Value = 276.5   # I117: 465 -> 276.5
$ws.Cells.Item(117, 10).Value = 1927.25   # J117: 1661.6364 -> 1927.25
$ws.Cells.Item(117, 11).Value = 829.5   # K117: 1395 -> 829.5
$ws.Cells.Item(117, 12).Value = 5781.75   # L117: 4984.9092 -> 5781.75
$ws.Cells.Item(117, 13).Value = 2612.5   # M117: 2047 -> 2612.5
$ws.Cells.Item(117, 14).Value = -12665.75   # N117: -11868.9092 -> -12665.75

$ws.Cells.Item(131, 8).Value = 667425.1   # H131: 1111815.5 -> 667425.1
$ws.Cells.Item(131, 10).Value = 715070.5   # J131: 1250743.8 -> 715070.5
$ws.Cells.Item(131, 12).Value = 2145211.5   # L131: 3752231.4 -> 2145211.5
$ws.Cells.Item(131, 14).Value = -2155291.5   # N131: -3762311.4 -> -2155291.5

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 22848.1   # H102: 34597.23 -> 22848.1
$ws.Cells.Item(102, 9).Value = 11532.308   # I102: 16315.111 -> 11532.308
$ws.Cells.Item(102, 10).Value = 43863.145   # J102: 75732 -> 43863.145
$ws.Cells.Item(102, 11).Value = 11532.308   # K102: 16315.111 -> 11532.308
$ws.Cells.Item(102, 12).Value = 43863.145   # L102: 75732 -> 43863.145
$ws.Cells.Item(102, 13).Value = -9910.308000000001   # M102: -14693.111 -> -9910.308000000001
$ws.Cells.Item(102, 14).Value = -47107.145   # N102: -78976 -> -47107.145

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(46, 8).Value = 2723.8333   # H46: 2879.6365 -> 2723.8333
$ws.Cells.Item(46, 10).Value = 3048.6   # J46: 3275.111 -> 3048.6
$ws.Cells.Item(46, 12).Value = 3048.6   # L46: 3275.111 -> 3048.6
$ws.Cells.Item(46, 14).Value = -3424.6   # N46: -3651.111 -> -3424.6

$ws.Cells.Item(50, 8).Value = 29800   # H50: 29790 -> 29800
$ws.Cells.Item(50, 10).Value = 29800   # J50: 29790 -> 29800
$ws.Cells.Item(50, 12).Value = 29800   # L50: 29790 -> 29800
$ws.Cells.Item(50, 14).Value = -31074   # N50: -31064 -> -31074

$ws.Cells.Item(68, 8).Value = 2036.9166   # H68: 2111 -> 2036.9166
$ws.Cells.Item(68, 10).Value = 2788.6   # J68: 2823.8333 -> 2788.6
$ws.Cells.Item(68, 12).Value = 2788.6   # L68: 2823.8333 -> 2788.6
$ws.Cells.Item(68, 14).Value = -4286.6   # N68: -4321.8333 -> -4286.6

$ws.Cells.Item(71, 8).Value = 2036.9166   # H71: 2111 -> 2036.9166
$ws.Cells.Item(71, 10).Value = 2788.6   # J71: 2823.8333 -> 2788.6
$ws.Cells.Item(71, 12).Value = 13943   # L71: 14119.1665 -> 13943
$ws.Cells.Item(71, 14).Value = -21431   # N71: -21607.1665 -> -21431

$ws.Cells.Item(82, 8).Value = 1510.6818   # H82: 1554.0476 -> 1510.6818
$ws.Cells.Item(82, 9).Value = 1156.8125   # I82: 1193.9333 -> 1156.8125
$ws.Cells.Item(82, 11).Value = 1156.8125   # K82: 1193.9333 -> 1156.8125
$ws.Cells.Item(82, 13).Value = -795.8125   # M82: -832.9332999999999 -> -795.8125

$ws.Cells.Item(85, 8).Value = 1510.6818   # H85: 1554.0476 -> 1510.6818
$ws.Cells.Item(85, 9).Value = 1156.8125   # I85: 1193.9333 -> 1156.8125
$ws.Cells.Item(85, 11).Value = 1156.8125   # K85: 1193.9333 -> 1156.8125
$ws.Cells.Item(85, 13).Value = 91.1875   # M85: 54.06670000000008 -> 91.1875

$ws.Cells.Item(122, 8).Value = 39218320   # H122: 2638.7837 -> 39218320
$ws.Cells.Item(122, 9).Value = 43480924   # I122: 2588.1482 -> 43480924
$ws.Cells.Item(122, 10).Value = 30305608   # J122: 2775.5 -> 30305608
$ws.Cells.Item(122, 11).Value = 130442772   # K122: 7764.444600000001 -> 130442772
$ws.Cells.Item(122, 12).Value = 90916824   # L122: 8326.5 -> 90916824
$ws.Cells.Item(122, 13).Value = -130440322   # M122: -5314.444600000001 -> -130440322
$ws.Cells.Item(122, 14).Value = -90921724   # N122: -13226.5 -> -90921724

$ws.Cells.Item(136, 8).Value = 477713.72   # H136: 456257.97 -> 477713.72
$ws.Cells.Item(136, 9).Value = 770135.25   # I136: 1251074.2 -> 770135.25
$ws.Cells.Item(136, 10).Value = 2528.75   # J136: 2077.2144 -> 2528.75
$ws.Cells.Item(136, 11).Value = 2310405.75   # K136: 3753222.6 -> 2310405.75
$ws.Cells.Item(136, 12).Value = 7586.25   # L136: 6231.6432 -> 7586.25
$ws.Cells.Item(136, 13).Value = -2307855.75   # M136: -3750672.6 -> -2307855.75
$ws.Cells.Item(136, 14).Value = -12686.25   # N136: -11331.6432 -> -12686.25

$ws.Cells.Item(138, 8).Value = 42433.332   # H138: 44171.43 -> 42433.332
$ws.Cells.Item(138, 10).Value = 42433.332   # J138: 44171.43 -> 42433.332
$ws.Cells.Item(138, 12).Value = 42433.332   # L138: 44171.43 -> 42433.332
$ws.Cells.Item(138, 14).Value = -52713.332   # N138: -54451.43 -> -52713.332

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(62, 8).Value = 4995.7144   # H62: 4952.222 -> 4995.7144
$ws.Cells.Item(62, 10).Value = 5333.3335   # J62: 5200 -> 5333.3335
$ws.Cells.Item(62, 12).Value = 5333.3335   # L62: 5200 -> 5333.3335
$ws.Cells.Item(62, 14).Value = -6581.3335   # N62: -6448 -> -6581.3335

$ws.Cells.Item(65, 8).Value = 4995.7144   # H65: 4952.222 -> 4995.7144
$ws.Cells.Item(65, 10).Value = 5333.3335   # J65: 5200 -> 5333.3335
$ws.Cells.Item(65, 12).Value = 26666.6675   # L65: 26000 -> 26666.6675
$ws.Cells.Item(65, 14).Value = -32906.6675   # N65: -32240 -> -32906.6675

$ws.Cells.Item(122, 8).Value = 5167.4443   # H122: 2656.6667 -> 5167.4443
$ws.Cells.Item(122, 9).Value = 1418.4286   # I122: 1322.2222 -> 1418.4286
$ws.Cells.Item(122, 10).Value = 7553.1816   # J122: 6660 -> 7553.1816
$ws.Cells.Item(122, 11).Value = 4255.2858   # K122: 3966.6666 -> 4255.2858
$ws.Cells.Item(122, 12).Value = 22659.5448   # L122: 19980 -> 22659.5448
$ws.Cells.Item(122, 13).Value = -1805.2858   # M122: -1516.6666 -> -1805.2858
$ws.Cells.Item(122, 14).Value = -27559.5448   # N122: -24880 -> -27559.5448

$ws.Cells.Item(133, 8).Value = 43855   # H133: 42236.668 -> 43855
$ws.Cells.Item(133, 10).Value = 43855   # J133: 42236.668 -> 43855
$ws.Cells.Item(133, 12).Value = 43855   # L133: 42236.668 -> 43855
$ws.Cells.Item(133, 14).Value = -53975   # N133: -52356.668 -> -53975
